$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Legislature" row (row 12) is a leftover/erroneous row for this
# de_dg_sw_upper parameter file. Delete the entire row so everything below
# it (Friant Water Authority, Friant-Kern Canal, etc.) shifts up by one,
# "adding rural communities interactions back in" by restoring the correct
# row alignment used by the other parameter files.
$ws.Rows.Item(12).Delete()

# Mirror the resulting selection/scroll position left behind by Excel after
# deleting the row (the row that slides up into row 12 ends up selected).
$ws.Range("A12:XFD12").Select()
$excel.ActiveWindow.ScrollRow = 10
